$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7399.5713
$ws.Range("I62").Value = 2359
$ws.Range("J62").Value = 20001
$ws.Range("K62").Value = 2359
$ws.Range("L62").Value = 20001
$ws.Range("M62").Value = -1735
$ws.Range("N62").Value = -21249
$ws.Range("H65").Value = 7399.5713
$ws.Range("I65").Value = 2359
$ws.Range("J65").Value = 20001
$ws.Range("K65").Value = 11795
$ws.Range("L65").Value = 100005
$ws.Range("M65").Value = -8675
$ws.Range("N65").Value = -106245
$ws.Range("H70").Value = 2176.647
$ws.Range("J70").Value = 2264.5
$ws.Range("L70").Value = 6793.5
$ws.Range("N70").Value = -7333.5
$ws.Range("H73").Value = 2176.647
$ws.Range("J73").Value = 2264.5
$ws.Range("L73").Value = 6793.5
$ws.Range("N73").Value = -8665.5
$ws.Range("H92").Value = 27778072
$ws.Range("I92").Value = 37037096
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 37037096
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = -37035848
$ws.Range("N92").Value = -3496
$ws.Range("H98").Value = 386592.38
$ws.Range("I98").Value = 400248.7
$ws.Range("J98").Value = 4216
$ws.Range("K98").Value = 400248.7
$ws.Range("L98").Value = 4216
$ws.Range("M98").Value = -398750.7
$ws.Range("N98").Value = -7212
$ws.Range("H122").Value = 386592.38
$ws.Range("I122").Value = 400248.7
$ws.Range("J122").Value = 4216
$ws.Range("K122").Value = 1200746.1
$ws.Range("L122").Value = 12648
$ws.Range("M122").Value = -1198296.1
$ws.Range("N122").Value = -17548
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H125").Value = 7474609.5
$ws.Range("I125").Value = 617.6923
$ws.Range("K125").Value = 5559.2307
$ws.Range("M125").Value = -3099.2307
$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H131").Value = 5903.933
$ws.Range("I131").Value = 1365.9
$ws.Range("J131").Value = 14980
$ws.Range("K131").Value = 4097.700000000001
$ws.Range("L131").Value = 44940
$ws.Range("M131").Value = 942.2999999999993
$ws.Range("N131").Value = -55020
$ws.Range("H138").Value = 2330.821
$ws.Range("I138").Value = 614.6667
$ws.Range("J138").Value = 2910.9297
$ws.Range("K138").Value = 1844.0001
$ws.Range("L138").Value = 8732.7891
$ws.Range("M138").Value = 3295.9999
$ws.Range("N138").Value = -19012.7891

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15854.027
$ws.Range("I32").Value = 2256.7188
$ws.Range("J32").Value = 112546
$ws.Range("K32").Value = 2256.7188
$ws.Range("L32").Value = 112546
$ws.Range("M32").Value = -1969.7188
$ws.Range("N32").Value = -113120
$ws.Range("H45").Value = 1418.75
$ws.Range("I45").Value = 1450
$ws.Range("J45").Value = 1387.5
$ws.Range("K45").Value = 1450
$ws.Range("L45").Value = 1387.5
$ws.Range("M45").Value = -1073
$ws.Range("N45").Value = -2141.5
$ws.Range("H74").Value = 12099.096
$ws.Range("I74").Value = 1816.875
$ws.Range("J74").Value = 45002.2
$ws.Range("K74").Value = 1816.875
$ws.Range("L74").Value = 45002.2
$ws.Range("M74").Value = -942.875
$ws.Range("N74").Value = -46750.2
$ws.Range("H77").Value = 12099.096
$ws.Range("I77").Value = 1816.875
$ws.Range("J77").Value = 45002.2
$ws.Range("K77").Value = 9084.375
$ws.Range("L77").Value = 225011
$ws.Range("M77").Value = -4716.375
$ws.Range("N77").Value = -233747
$ws.Range("H122").Value = 1566.6364
$ws.Range("I122").Value = 1249.75
$ws.Range("J122").Value = 2121.1875
$ws.Range("K122").Value = 3749.25
$ws.Range("L122").Value = 6363.5625
$ws.Range("M122").Value = -1299.25
$ws.Range("N122").Value = -11263.5625
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4201.027
$ws.Range("I86").Value = 1554.5416
$ws.Range("J86").Value = 9086.846
$ws.Range("K86").Value = 1554.5416
$ws.Range("L86").Value = 9086.846
$ws.Range("M86").Value = -431.5416
$ws.Range("N86").Value = -11332.846
$ws.Range("H89").Value = 4201.027
$ws.Range("I89").Value = 1554.5416
$ws.Range("J89").Value = 9086.846
$ws.Range("K89").Value = 7772.708000000001
$ws.Range("L89").Value = 45434.23
$ws.Range("M89").Value = -2156.708000000001
$ws.Range("N89").Value = -56666.23
$ws.Range("H94").Value = 1093.3478
$ws.Range("I94").Value = 1145.1666
$ws.Range("J94").Value = 906.8
$ws.Range("K94").Value = 1145.1666
$ws.Range("L94").Value = 906.8
$ws.Range("M94").Value = -694.1666
$ws.Range("N94").Value = -1808.8
$ws.Range("H134").Value = 4511.591
$ws.Range("I134").Value = 2944.3
$ws.Range("J134").Value = 5817.6665
$ws.Range("K134").Value = 8832.900000000001
$ws.Range("L134").Value = 17452.9995
$ws.Range("M134").Value = -6297.900000000001
$ws.Range("N134").Value = -22522.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8929845
$ws.Range("I99").Value = 12501202
$ws.Range("J99").Value = 1450
$ws.Range("K99").Value = 12501202
$ws.Range("L99").Value = 1450
$ws.Range("M99").Value = -12499704
$ws.Range("N99").Value = -4446
$ws.Range("H126").Value = 8929845
$ws.Range("I126").Value = 12501202
$ws.Range("J126").Value = 1450
$ws.Range("K126").Value = 37503606
$ws.Range("L126").Value = 4350
$ws.Range("M126").Value = -37501136
$ws.Range("N126").Value = -9290
$ws.Range("H127").Value = 35000
$ws.Range("J127").Value = 35000
$ws.Range("L127").Value = 35000
$ws.Range("N127").Value = -44920
$ws.Range("H134").Value = 2930.3447
$ws.Range("I134").Value = 1644.3684
$ws.Range("J134").Value = 5373.7
$ws.Range("K134").Value = 4933.1052
$ws.Range("L134").Value = 16121.1
$ws.Range("M134").Value = -2398.1052
$ws.Range("N134").Value = -21191.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 1500
$ws.Range("I116").Value = 1166.6666
$ws.Range("K116").Value = 3499.9998
$ws.Range("M116").Value = -57.99980000000005
$ws.Range("H117").Value = 1140
$ws.Range("H131").Value = 1354.0769
$ws.Range("I131").Value = 320.83334
$ws.Range("J131").Value = 1588.0189
$ws.Range("K131").Value = 962.5000200000001
$ws.Range("L131").Value = 4764.0567
$ws.Range("M131").Value = 4077.49998
$ws.Range("N131").Value = -14844.0567

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 44
$ws.Range("I2").Value = 31.5
$ws.Range("J2").Value = 54
$ws.Range("K2").Value = 31.5
$ws.Range("L2").Value = 54
$ws.Range("M2").Value = 81.5
$ws.Range("N2").Value = -280
$ws.Range("H70").Value = 6981.579
$ws.Range("I70").Value = 7511.5386
$ws.Range("J70").Value = 5833.3335
$ws.Range("K70").Value = 7511.5386
$ws.Range("L70").Value = 5833.3335
$ws.Range("M70").Value = -7241.5386
$ws.Range("N70").Value = -6373.3335
$ws.Range("H73").Value = 6981.579
$ws.Range("I73").Value = 7511.5386
$ws.Range("J73").Value = 5833.3335
$ws.Range("K73").Value = 7511.5386
$ws.Range("L73").Value = 5833.3335
$ws.Range("M73").Value = -6575.5386
$ws.Range("N73").Value = -7705.3335
$ws.Range("H123").Value = 9728.091
$ws.Range("J123").Value = 9728.091
$ws.Range("L123").Value = 9728.091
$ws.Range("N123").Value = -14628.091
$ws.Range("H132").Value = 4245.6055
$ws.Range("I132").Value = 4201.56
$ws.Range("J132").Value = 4330.3076
$ws.Range("K132").Value = 12604.68
$ws.Range("L132").Value = 12990.9228
$ws.Range("M132").Value = -10074.68
$ws.Range("N132").Value = -18050.9228
$ws.Range("H133").Value = 16794.445
$ws.Range("J133").Value = 16794.445
$ws.Range("L133").Value = 16794.445
$ws.Range("N133").Value = -26914.445

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1258.3334
$ws.Range("I22").Value = 450
$ws.Range("J22").Value = 1662.5
$ws.Range("K22").Value = 450
$ws.Range("L22").Value = 1662.5
$ws.Range("M22").Value = -155
$ws.Range("N22").Value = -2252.5
$ws.Range("H27").Value = 1258.3334
$ws.Range("I27").Value = 450
$ws.Range("J27").Value = 1662.5
$ws.Range("K27").Value = 450
$ws.Range("L27").Value = 1662.5
$ws.Range("M27").Value = -343
$ws.Range("N27").Value = -1876.5
$ws.Range("H132").Value = 3278.756
$ws.Range("I132").Value = 2689.1177
$ws.Range("J132").Value = 6142.7144
$ws.Range("K132").Value = 8067.353099999999
$ws.Range("L132").Value = 18428.1432
$ws.Range("M132").Value = -5537.353099999999
$ws.Range("N132").Value = -23488.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16689069
$ws.Range("I62").Value = 25028076
$ws.Range("J62").Value = 11055
$ws.Range("K62").Value = 25028076
$ws.Range("L62").Value = 11055
$ws.Range("M62").Value = -25027452
$ws.Range("N62").Value = -12303
$ws.Range("H65").Value = 16689069
$ws.Range("I65").Value = 25028076
$ws.Range("J65").Value = 11055
$ws.Range("K65").Value = 125140380
$ws.Range("L65").Value = 55275
$ws.Range("M65").Value = -125137260
$ws.Range("N65").Value = -61515
$ws.Range("H128").Value = 49900
$ws.Range("J128").Value = 49900
$ws.Range("L128").Value = 49900
$ws.Range("N128").Value = -59860
$ws.Range("H136").Value = 8359854.5
$ws.Range("I136").Value = 10132325
$ws.Range("J136").Value = 3921.5715
$ws.Range("K136").Value = 30396975
$ws.Range("L136").Value = 11764.7145
$ws.Range("M136").Value = -30394425
$ws.Range("N136").Value = -16864.7145
